{"js": "// Update each \"AA\u00d7BB=CCCC\" multiplication-answer cell in the table to the\n// new equation/result given in the commit's target revision. Every old\n// value appears exactly once in the document, so a scoped search +\n// insertText(\"Replace\") per pair is unambiguous and order-independent.\nconst replacements = [\n  [\"54\u00d742=2268\", \"30\u00d752=1560\"],\n  [\"28\u00d791=2548\", \"76\u00d782=6232\"],\n  [\"74\u00d752=3848\", \"27\u00d784=2268\"],\n  [\"52\u00d736=1872\", \"63\u00d727=1701\"],\n  [\"97\u00d779=7663\", \"21\u00d754=1134\"],\n  [\"83\u00d750=4150\", \"39\u00d744=1716\"],\n  [\"42\u00d746=1932\", \"16\u00d786=1376\"],\n  [\"98\u00d780=7840\", \"98\u00d772=7056\"],\n  [\"85\u00d786=7310\", \"33\u00d714=462\"],\n  [\"63\u00d759=3717\", \"61\u00d730=1830\"],\n  [\"82\u00d722=1804\", \"38\u00d738=1444\"],\n  [\"76\u00d790=6840\", \"67\u00d771=4757\"],\n  [\"28\u00d753=1484\", \"56\u00d722=1232\"],\n  [\"24\u00d756=1344\", \"77\u00d780=6160\"],\n  [\"47\u00d722=1034\", \"89\u00d760=5340\"],\n  [\"16\u00d798=1568\", \"27\u00d721=567\"],\n  [\"51\u00d732=1632\", \"80\u00d768=5440\"],\n  [\"61\u00d784=5124\", \"67\u00d770=4690\"],\n  [\"45\u00d760=2700\", \"37\u00d727=999\"],\n  [\"86\u00d720=1720\", \"42\u00d775=3150\"],\n  [\"40\u00d746=1840\", \"75\u00d754=4050\"],\n  [\"51\u00d729=1479\", \"19\u00d765=1235\"],\n  [\"42\u00d792=3864\", \"65\u00d717=1105\"],\n  [\"91\u00d757=5187\", \"39\u00d716=624\"],\n  [\"47\u00d768=3196\", \"85\u00d735=2975\"],\n];\n\nconst body = context.document.body;\n\nfor (const [oldText, newText] of replacements) {\n  const results = body.search(oldText, { matchCase: true, matchWholeWord: false });\n  results.load(\"items\");\n  await context.sync();\n\n  for (const range of results.items) {\n    range.insertText(newText, \"Replace\");\n  }\n  await context.sync();\n}\n", "ps1": "# Replace each two-digit x two-digit multiplication \"equation\" cell\n# text with its updated value. Each old value is unique in the\n# document, so a plain Find/Replace (ReplaceAll) per pair is safe and\n# unambiguous. Scoped to the whole document body via $d.Content.\n$d = $word.ActiveDocument\n\n$replacements = @(\n    @(\"54\u00d742=2268\", \"30\u00d752=1560\"),\n    @(\"28\u00d791=2548\", \"76\u00d782=6232\"),\n    @(\"74\u00d752=3848\", \"27\u00d784=2268\"),\n    @(\"52\u00d736=1872\", \"63\u00d727=1701\"),\n    @(\"97\u00d779=7663\", \"21\u00d754=1134\"),\n    @(\"83\u00d750=4150\", \"39\u00d744=1716\"),\n    @(\"42\u00d746=1932\", \"16\u00d786=1376\"),\n    @(\"98\u00d780=7840\", \"98\u00d772=7056\"),\n    @(\"85\u00d786=7310\", \"33\u00d714=462\"),\n    @(\"63\u00d759=3717\", \"61\u00d730=1830\"),\n    @(\"82\u00d722=1804\", \"38\u00d738=1444\"),\n    @(\"76\u00d790=6840\", \"67\u00d771=4757\"),\n    @(\"28\u00d753=1484\", \"56\u00d722=1232\"),\n    @(\"24\u00d756=1344\", \"77\u00d780=6160\"),\n    @(\"47\u00d722=1034\", \"89\u00d760=5340\"),\n    @(\"16\u00d798=1568\", \"27\u00d721=567\"),\n    @(\"51\u00d732=1632\", \"80\u00d768=5440\"),\n    @(\"61\u00d784=5124\", \"67\u00d770=4690\"),\n    @(\"45\u00d760=2700\", \"37\u00d727=999\"),\n    @(\"86\u00d720=1720\", \"42\u00d775=3150\"),\n    @(\"40\u00d746=1840\", \"75\u00d754=4050\"),\n    @(\"51\u00d729=1479\", \"19\u00d765=1235\"),\n    @(\"42\u00d792=3864\", \"65\u00d717=1105\"),\n    @(\"91\u00d757=5187\", \"39\u00d716=624\"),\n    @(\"47\u00d768=3196\", \"85\u00d735=2975\")\n)\n\nforeach ($pair in $replacements) {\n    $oldText = $pair[0]\n    $newText = $pair[1]\n    $find = $d.Content.Find\n    $find.ClearFormatting()\n    $find.Text = $oldText\n    $find.Replacement.ClearFormatting()\n    $find.Replacement.Text = $newText\n    $find.Execute(\n        $oldText,        # FindText\n        $true,           # MatchCase\n        $false,          # MatchWholeWord\n        $false,          # MatchWildcards\n        $false,          # MatchSoundsLike\n        $false,          # MatchAllWordForms\n        $true,           # Forward\n        1,               # Wrap (wdFindContinue)\n        $false,          # Format\n        $newText,        # ReplaceWith\n        2                # Replace (wdReplaceAll)\n    ) | Out-Null\n}\n"}
